$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates ---
$ws.Range("A8").Value = "Volume 32   Number  45"
$ws.Range("C9").Value = "Report Covering the Week  11/3/2025  Through  11/9/2025"

# --- Weekly crime statistics table (rows 14-33) ---
# Row 14: Murder
$ws.Range("A14").Value = "Murder"
$ws.Range("C14").Value = "'0"
$ws.Range("D14").Value = "'0"
$ws.Range("E14").Value = "***.*"
$ws.Range("F14").Value = 2
$ws.Range("G14").Value = 6
$ws.Range("H14").Value = -66.666666666666
$ws.Range("I14").Value = 34
$ws.Range("J14").Value = 46
$ws.Range("K14").Value = -26.086956521739
$ws.Range("L14").Value = -37.037037037037
$ws.Range("M14").Value = -54.666666666666
$ws.Range("N14").Value = -84.40366972477

# Row 15: Rape
$ws.Range("A15").Value = "Rape"
$ws.Range("C15").Value = 4
$ws.Range("D15").Value = 3
$ws.Range("E15").Value = 33.333333333333
$ws.Range("F15").Value = 21
$ws.Range("G15").Value = 16
$ws.Range("H15").Value = 31.25
$ws.Range("I15").Value = 227
$ws.Range("J15").Value = 199
$ws.Range("K15").Value = 14.070351758794
$ws.Range("L15").Value = 21.39037433155
$ws.Range("M15").Value = 47.402597402597
$ws.Range("N15").Value = -53.955375253549

# Row 16: Robbery
$ws.Range("A16").Value = "Robbery"
$ws.Range("C16").Value = 44
$ws.Range("D16").Value = 38
$ws.Range("E16").Value = 15.78947368421
$ws.Range("F16").Value = 142
$ws.Range("G16").Value = 131
$ws.Range("H16").Value = 8.396946564885
$ws.Range("I16").Value = 1428
$ws.Range("J16").Value = 1497
$ws.Range("K16").Value = -4.609218436873
$ws.Range("L16").Value = -6.11439842209
$ws.Range("M16").Value = -43.084894380231
$ws.Range("N16").Value = -88.136578881781

# Row 17: Fel. Assault
$ws.Range("A17").Value = "Fel. Assault"
$ws.Range("C17").Value = 84
$ws.Range("D17").Value = 80
$ws.Range("E17").Value = 5
$ws.Range("F17").Value = 322
$ws.Range("G17").Value = 292
$ws.Range("H17").Value = 10.273972602739
$ws.Range("I17").Value = 3405
$ws.Range("J17").Value = 3234
$ws.Range("K17").Value = 5.287569573283
$ws.Range("L17").Value = 11.383709519136
$ws.Range("M17").Value = 60.764872521246
$ws.Range("N17").Value = -41.665238992633

# Row 18: Burglary
$ws.Range("A18").Value = "Burglary"
$ws.Range("C18").Value = 20
$ws.Range("D18").Value = 32
$ws.Range("E18").Value = -37.5
$ws.Range("F18").Value = 96
$ws.Range("G18").Value = 135
$ws.Range("H18").Value = -28.888888888888
$ws.Range("I18").Value = 1250
$ws.Range("J18").Value = 1330
$ws.Range("K18").Value = -6.015037593984
$ws.Range("L18").Value = -15.76819407008
$ws.Range("M18").Value = -58.388814913448
$ws.Range("N18").Value = -92.168901140208

# Row 19: Gr. Larceny
$ws.Range("A19").Value = "Gr. Larceny"
$ws.Range("C19").Value = 116
$ws.Range("D19").Value = 121
$ws.Range("E19").Value = -4.132231404958
$ws.Range("F19").Value = 404
$ws.Range("G19").Value = 442
$ws.Range("H19").Value = -8.597285067873
$ws.Range("I19").Value = 4883
$ws.Range("J19").Value = 5056
$ws.Range("K19").Value = -3.421677215189
$ws.Range("L19").Value = -13.621086148947
$ws.Range("M19").Value = 3.343915343915
$ws.Range("N19").Value = -36.128188358404

# Row 20: G.L.A.
$ws.Range("A20").Value = "G.L.A."
$ws.Range("C20").Value = 32
$ws.Range("D20").Value = 34
$ws.Range("E20").Value = -5.882352941176
$ws.Range("F20").Value = 119
$ws.Range("G20").Value = 156
$ws.Range("H20").Value = -23.717948717948
$ws.Range("I20").Value = 1570
$ws.Range("J20").Value = 1759
$ws.Range("K20").Value = -10.744741330301
$ws.Range("L20").Value = -2.665840049597
$ws.Range("M20").Value = -6.21266427718
$ws.Range("N20").Value = -92.16293116358

# Row 21: TOTAL
$ws.Range("A21").Value = "TOTAL"
$ws.Range("C21").Value = 300
$ws.Range("D21").Value = 308
$ws.Range("E21").Value = -2.597402597402
$ws.Range("F21").Value = 1106
$ws.Range("G21").Value = 1178
$ws.Range("H21").Value = -6.112054329371
$ws.Range("I21").Value = 12797
$ws.Range("J21").Value = 13121
$ws.Range("K21").Value = -2.469323984452
$ws.Range("L21").Value = -5.689439162797
$ws.Range("M21").Value = -10.253173434322
$ws.Range("N21").Value = -79.43431096826

# Row 22: Transit
$ws.Range("A22").Value = "Transit"
$ws.Range("C22").Value = 5
$ws.Range("D22").Value = 3
$ws.Range("E22").Value = 66.666666666666
$ws.Range("F22").Value = 24
$ws.Range("G22").Value = 15
$ws.Range("H22").Value = 60
$ws.Range("I22").Value = 198
$ws.Range("J22").Value = 174
$ws.Range("K22").Value = 13.793103448275
$ws.Range("L22").Value = 13.793103448275
$ws.Range("M22").Value = -18.852459016393
$ws.Range("N22").Value = "***.*"

# Row 23: Housing
$ws.Range("A23").Value = "Housing"
$ws.Range("C23").Value = 7
$ws.Range("D23").Value = 15
$ws.Range("E23").Value = -53.333333333333
$ws.Range("F23").Value = 32
$ws.Range("G23").Value = 48
$ws.Range("H23").Value = -33.333333333333
$ws.Range("I23").Value = 416
$ws.Range("J23").Value = 435
$ws.Range("K23").Value = -4.367816091954
$ws.Range("L23").Value = -8.370044052863
$ws.Range("M23").Value = 41.016949152542
$ws.Range("N23").Value = "***.*"

# Row 24: Petit Larceny
$ws.Range("A24").Value = "Petit Larceny"
$ws.Range("C24").Value = 262
$ws.Range("D24").Value = 330
$ws.Range("E24").Value = -20.60606060606
$ws.Range("F24").Value = 1002
$ws.Range("G24").Value = 1155
$ws.Range("H24").Value = -13.246753246753
$ws.Range("I24").Value = 11594
$ws.Range("J24").Value = 12876
$ws.Range("K24").Value = -9.95650823237
$ws.Range("L24").Value = -15.149297423887
$ws.Range("M24").Value = 8.314648729446
$ws.Range("N24").Value = "***.*"

# Row 25: Retail Theft
$ws.Range("A25").Value = "Retail Theft"
$ws.Range("C25").Value = 106
$ws.Range("D25").Value = 149
$ws.Range("E25").Value = -28.859060402684
$ws.Range("F25").Value = 406
$ws.Range("G25").Value = 524
$ws.Range("H25").Value = -22.519083969465
$ws.Range("I25").Value = 4669
$ws.Range("J25").Value = 6137
$ws.Range("K25").Value = -23.920482320352
$ws.Range("L25").Value = -29.428657799274
$ws.Range("M25").Value = "***.*"
$ws.Range("N25").Value = "***.*"

# Row 26: Misd. Assault
$ws.Range("A26").Value = "Misd. Assault"
$ws.Range("C26").Value = 140
$ws.Range("D26").Value = 133
$ws.Range("E26").Value = 5.263157894736
$ws.Range("F26").Value = 541
$ws.Range("G26").Value = 561
$ws.Range("H26").Value = -3.565062388591
$ws.Range("I26").Value = 5565
$ws.Range("J26").Value = 5644
$ws.Range("K26").Value = -1.399716513111
$ws.Range("L26").Value = 9.871668311944
$ws.Range("M26").Value = -4.839261285909
$ws.Range("N26").Value = "***.*"

# Row 27: UCR Rape*
$ws.Range("A27").Value = "UCR Rape*"
$ws.Range("C27").Value = 4
$ws.Range("D27").Value = 3
$ws.Range("E27").Value = 33.333333333333
$ws.Range("F27").Value = 24
$ws.Range("G27").Value = 18
$ws.Range("H27").Value = 33.333333333333
$ws.Range("I27").Value = 264
$ws.Range("J27").Value = 287
$ws.Range("K27").Value = -8.013937282229
$ws.Range("L27").Value = -2.583025830258
$ws.Range("M27").Value = "***.*"
$ws.Range("N27").Value = "***.*"

# Row 28: Other Sex Crimes
$ws.Range("A28").Value = "Other Sex Crimes"
$ws.Range("C28").Value = 16
$ws.Range("D28").Value = 8
$ws.Range("E28").Value = 100
$ws.Range("F28").Value = 59
$ws.Range("G28").Value = 44
$ws.Range("H28").Value = 34.090909090909
$ws.Range("I28").Value = 588
$ws.Range("J28").Value = 594
$ws.Range("K28").Value = -1.010101010101
$ws.Range("L28").Value = 5.945945945945
$ws.Range("M28").Value = "***.*"
$ws.Range("N28").Value = "***.*"

# Row 29: Shooting Vic.
$ws.Range("A29").Value = "Shooting Vic."
$ws.Range("C29").Value = 4
$ws.Range("D29").Value = "'0"
$ws.Range("E29").Value = "***.*"
$ws.Range("F29").Value = 10
$ws.Range("G29").Value = 5
$ws.Range("H29").Value = 100
$ws.Range("I29").Value = 116
$ws.Range("J29").Value = 92
$ws.Range("K29").Value = 26.086956521739
$ws.Range("L29").Value = -7.936507936507
$ws.Range("M29").Value = -52.459016393442
$ws.Range("N29").Value = -83.164005805515

# Row 30: Shooting Inc.
$ws.Range("A30").Value = "Shooting Inc."
$ws.Range("C30").Value = 4
$ws.Range("D30").Value = "'0"
$ws.Range("E30").Value = "***.*"
$ws.Range("F30").Value = 10
$ws.Range("G30").Value = 5
$ws.Range("H30").Value = 100
$ws.Range("I30").Value = 89
$ws.Range("J30").Value = 79
$ws.Range("K30").Value = 12.658227848101
$ws.Range("L30").Value = -18.348623853211
$ws.Range("M30").Value = -55.940594059405
$ws.Range("N30").Value = -85.215946843853

# Row 31: Hate Crimes
$ws.Range("A31").Value = "Hate Crimes"
$ws.Range("C31").Value = 2
$ws.Range("D31").Value = 3
$ws.Range("E31").Value = -33.333333333333
$ws.Range("F31").Value = 6
$ws.Range("G31").Value = 11
$ws.Range("H31").Value = -45.454545454545
$ws.Range("I31").Value = 103
$ws.Range("J31").Value = 139
$ws.Range("K31").Value = -25.899280575539
$ws.Range("L31").Value = 9.574468085106
$ws.Range("M31").Value = "***.*"
$ws.Range("N31").Value = "***.*"

# Row 33: Traffic Fatalities
$ws.Range("A33").Value = "Traffic Fatalities"
$ws.Range("C33").Value = "'0"
$ws.Range("D33").Value = "'0"
$ws.Range("E33").Value = "***.*"
$ws.Range("F33").Value = 4
$ws.Range("G33").Value = 3
$ws.Range("H33").Value = 33.333333333333
$ws.Range("I33").Value = 37
$ws.Range("J33").Value = 36
$ws.Range("K33").Value = 2.777777777777
$ws.Range("L33").Value = -7.5
$ws.Range("M33").Value = "***.*"
$ws.Range("N33").Value = "***.*"

# --- Column E width fix (narrower, since % Chg no longer needs extra width) ---
$ws.Range("E1").EntireColumn.ColumnWidth = 6.168446